$wb = $excel.ActiveWorkbook

# Update the time_taken values on the "data" sheet (F2:F8)
$data = $wb.Worksheets.Item("data")
$data.Range("F2").Value = "2021-10-05 14:20:12.705370"
$data.Range("F3").Value = "2021-10-05 14:20:12.705378"
$data.Range("F4").Value = "2021-10-05 14:20:12.705381"
$data.Range("F5").Value = "2021-10-05 14:20:12.705384"
$data.Range("F6").Value = "2021-10-05 14:20:12.705386"
$data.Range("F7").Value = "2021-10-05 14:20:12.705389"
$data.Range("F8").Value = "2021-10-05 14:20:12.705392"

# Add the new "metadata" sheet directly after "data"
$data = $wb.Worksheets.Item("data")
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Re-fetch fresh handles (sheet indices shifted after Add)
$meta = $wb.Worksheets.Item("metadata")
$data = $wb.Worksheets.Item("data")

# Match the "data" sheet's page margins (0.75in/0.75in/1in/1in/0.5in/0.5in -> points)
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Match the bold/centered/bordered header style used on the "data" sheet
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Familial disseminated superficial actinic porokeratosis"
$meta.Range("C2").Value = 110
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.1"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2017-11-05T02:37:20.048842Z"
$meta.Range("F2").Value = "2021-10-05 14:20:12.701720"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/110/?format=json"
